# fix response and password generator while upload data siswa
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New names introduced in this order: Asti, Damah, Anis, then gender P
$ws.Cells.Item(3, 2).Value = "Asti"
$ws.Cells.Item(2, 2).Value = "Damah"
$ws.Cells.Item(4, 2).Value = "Anis"

$ws.Cells.Item(2, 4).Value = "P"
$ws.Cells.Item(3, 4).Value = "P"
$ws.Cells.Item(4, 4).Value = "P"

$ws.Cells.Item(2, 5).Value = "profile.png"
$ws.Cells.Item(3, 5).Value = "profile.png"
$ws.Cells.Item(4, 5).Value = "profile.png"

$ws.Cells.Item(2, 1).Value = 990229
$ws.Cells.Item(3, 1).Value = 921200
$ws.Cells.Item(4, 1).Value = 231213

$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(4, 3).ClearContents()
